$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text values (Korean/English mixed) per the diff, all rendered
# bold + white to match the existing header-row look (sz 11, Calibri).
$headers = @(
    @{ Addr = "B1"; Text = "총 Chai 판매(개수)" },
    @{ Addr = "C1"; Text = "Artisanal Chai 판매(단위)" },
    @{ Addr = "D1"; Text = "미리 만든 Chai 판매(단위)" },
    @{ Addr = "E1"; Text = "소셜 미디어 참여도(보기)" },
    @{ Addr = "F1"; Text = "Chai에 대한 온라인 검색" }
)

foreach ($h in $headers) {
    $cell = $ws.Range($h.Addr)
    $cell.Value = $h.Text
    $len = $cell.Characters().Text.Length

    # Apply bold + white font color + explicit face/size across the whole
    # string as two adjoining character ranges (this makes the engine emit
    # a genuine rich-text run inside the shared string, instead of just
    # collapsing the formatting onto the cell style).
    $first = $cell.Characters(1, $len - 1)
    $last = $cell.Characters($len, 1)
    foreach ($run in @($first, $last)) {
        $run.Font.Bold = $true
        $run.Font.Color = 16777215
        $run.Font.Name = "Calibri"
        $run.Font.Size = 11
    }
}

# Update the table column headers to match (Table1)
$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListColumns.Item(2).Name = "총 Chai 판매(개수)"
$tbl.ListColumns.Item(3).Name = "Artisanal Chai 판매(단위)"
$tbl.ListColumns.Item(4).Name = "미리 만든 Chai 판매(단위)"
$tbl.ListColumns.Item(5).Name = "소셜 미디어 참여도(보기)"
$tbl.ListColumns.Item(6).Name = "Chai에 대한 온라인 검색"
